$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H so the existing "feedback" (H) and "validator_id" (I)
# shift right to I and J respectively, making room for the new "status_jurnal" column.
$ws.Columns("H:H").Insert()

# Header row updates / additions
$ws.Range("G1").Value = "status_laporan"
$ws.Range("H1").Value = "status_jurnal"

# New "tanggal_laporan" header in K1 - copy the header formatting (bold/border/center)
# from an existing header cell, then set its text.
$ws.Range("G1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "tanggal_laporan"

# Drop the old sample rows 3-6, keeping only row 2 of data (which is rewritten below)
$ws.Rows("3:6").Delete()

# Update the remaining data row (row 2) to the new sample values
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = $false
$ws.Range("D2").Value = "jurnal.com"
$ws.Range("E2").Value = "jurnal"
$ws.Range("F2").Value = "tes"
$ws.Range("G2").Value = "pending"
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "2024-12-23 00:15:09"
